$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: split "Lab 7: Variable Selection" into label (C25) + moved text (H25) ---
# (Done first so new shared strings are appended in the same order as the target file)
$ws.Range("C25").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("H25").Value = "Lab 7: Variable Selection"
$ws.Range("C25").Value = "Lab 7: "

# --- Row 18: add Lecture 11 reading/slides links (D18, E18) ---
# Copy formatting from C18 (style index 2) into D18 and E18, then set values
$ws.Range("C18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("D18").Value = "11-reading.html"
$ws.Range("E18").Value = "11-bayes-regression"

# --- Update the active selection on the sheet to E18 ---
$ws.Range("E18").Select() | Out-Null

$excel.CutCopyMode = 0
